# Updated several exercises relating to SWC 2.semester
# ArithmeticConversion.pptx edits:
#  1) Delete the "int age = age + 2" example slide (position 2).
#  2) Refresh the cached "datetimeFigureOut" date field text
#     (01-02-2018 -> 11-02-2018) on the slide master and on every
#     slide layout.

$p = $ppt.ActivePresentation

# --- 1. Delete the obsolete "age + 2" slide -------------------------------
for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $slide = $p.Slides.Item($i)
    $text = ""
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shp = $slide.Shapes.Item($j)
        if ($shp.HasTextFrame -eq -1) {
            $text = $text + $shp.TextFrame.TextRange.Text
        }
    }
    if ($text -like "*age + 2*") {
        $slide.Delete()
    }
}

# --- 2. Update the cached date field everywhere ---------------------------
$master = $p.SlideMaster

for ($j = 1; $j -le $master.Shapes.Count; $j++) {
    $shp = $master.Shapes.Item($j)
    if ($shp.HasTextFrame -eq -1) {
        if ($shp.TextFrame.TextRange.Text -eq "01-02-2018") {
            $shp.TextFrame.TextRange.Text = "11-02-2018"
        }
    }
}

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $shp = $layout.Shapes.Item($j)
        if ($shp.HasTextFrame -eq -1) {
            if ($shp.TextFrame.TextRange.Text -eq "01-02-2018") {
                $shp.TextFrame.TextRange.Text = "11-02-2018"
            }
        }
    }
}
